$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 0, 1, 0.006966345516345562),
    @(2, 1, 1, 0.01220146782455322),
    @(3, 2, 1, 0.0003567146545742084),
    @(5, 4, 1, [double]"5.819229140002768e-06"),
    @(6, 5, 1, [double]"4.572818349535268e-07"),
    @(7, 5, 2, [double]"2.592921699129391e-06"),
    @(8, 7, 1, [double]"3.485003710057555e-07"),
    @(10, 8, 2, [double]"1.782396432670197e-07")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
